# DOMA-1872: remove multi-tariff values for non-electricity meters from
# meter-import-example.xlsx
#
# Rows (1-based, header in row 1) whose meter type (column D) is NOT "ЭЛ"
# (electricity) get their tariff count (column F) reset to 1 and their
# second/third readings (columns H/I) cleared, since non-electricity
# meters are single-tariff. Electricity rows (6 & 7, "ЭЛ") stay untouched.
# Also fixes the counter number for the gas meter row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nonElectricityRows = @(2, 3, 4, 5, 8, 9, 10)

foreach ($row in $nonElectricityRows) {
    $ws.Range("F$row").Value = 1
    $ws.Range("H$row").ClearContents()
    $ws.Range("I$row").ClearContents()
}

# Gas meter counter number correction (row 10, column E: "№ счетчика")
$ws.Range("E10").Value = 33
